$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 37-39: Nu/Inter/XP bank ticker mappings with P/VP and Div.Yield
# figures. Shared strings must be introduced in this exact order so the
# rebuilt sst table lines up with the target workbook (328..336):
#   INBR, ROXO, XPBR, 31,8, 9,97, 16,0, 0,97%, 1,07%, 0,0%

$ws.Range("A37").Value = "INBR"
$ws.Range("A38").Value = "ROXO"
$ws.Range("A39").Value = "XPBR"

$ws.Range("C38").Value = "31,8"
$ws.Range("C39").Value = "9,97"
$ws.Range("C37").Value = "16,0"

$ws.Range("F37").Value = "0,97%"
$ws.Range("F37").NumberFormat = "0.00%"

$ws.Range("F39").Value = "1,07%"
$ws.Range("F39").NumberFormat = "0.00%"

$ws.Range("F38").Value = "0,0%"
$ws.Range("F38").NumberFormat = "0%"

# Match the saved selection/scroll state from the edit.
$ws.Range("B37:B39").Select()
